$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.833.12'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.40%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.449.17'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.96%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.78'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.28'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -4.64%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("E8").Value = '  -3.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.449.88'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.85%  '

$ws.Range("E10").Value = '  -3.97%  '

$ws.Range("E11").Value = '  -1.35%  '

$ws.Range("E12").Value = '  -3.42%  '

$ws.Range("E13").Value = '  -5.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.893.18'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.01'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -6.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.773.12'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000166'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -6.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.441.67'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.87'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -9.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.32'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -9.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '347.88'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.85%  '

$ws.Range("E22").Value = '  -4.92%  '

$ws.Range("E23").Value = '  -2.29%  '

$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '68.37'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.15'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -9.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.78'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -6.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.97'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -10.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.996'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -43.69%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.573.38'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.67%  '

$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '503.64'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -5.96%  '

$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0₃0882'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -9.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.54'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -9.56%  '

$ws.Range("E34").Value = '  -6.72%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.20'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -8.62%  '

$ws.Range("E36").Value = '  +0.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '158.07'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.35%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.112'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -13.41%  '

$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.05'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -6.53%  '

$ws.Range("E41").Value = '  -9.56%  '

$ws.Range("E42").Value = '  +0.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.65'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -7.97%  '

$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.321'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -8.00%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.70'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -8.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.33'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -7.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '38.52'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '139.81'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -6.10%  '

$ws.Range("E49").Value = '  -8.98%  '

$ws.Range("E50").Value = '  -9.83%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0724'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.46%  '
